# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 45178 (2023-09-09) to 45179 (2023-09-10).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 45179
}
